{"js": "// Fix missing expert testimony achievement in comprehensive resumes.\n// Inserts four new bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n// section, right after the \"Platform impact...\" bullet and right before the\n// \"TECHNICAL SKILLS\" heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph: the last bullet under \"KEY ACHIEVEMENTS AND\n// IMPACT\" (\"\u2022 Platform impact: Built redistricting system serving 12,847\n// analysts across 89 organizations\").\nconst anchorText = \"Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\";\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate the 'Platform impact' achievement bullet paragraph.\");\n}\n\n// The four bullets to add, in document order.\nconst newBullets = [\n  \"\\u2022 Real-time collaboration at national scale\",\n  \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n  \"\\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n  \"\\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"\n];\n\n// Insert each bullet as its own new paragraph, immediately after the anchor,\n// chaining so they appear in the same order as in the diff.\nlet insertAfter = anchorParagraph;\nfor (const bulletText of newBullets) {\n  insertAfter = insertAfter.insertParagraph(bulletText, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Fix missing expert testimony achievement in comprehensive resumes.\n# Inserts four new bullet paragraphs into the \"KEY ACHIEVEMENTS AND IMPACT\"\n# section, right after the \"Platform impact...\" bullet and right before the\n# \"TECHNICAL SKILLS\" heading.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: the last bullet under \"KEY ACHIEVEMENTS AND\n# IMPACT\" (\"\u2022 Platform impact: Built redistricting system serving 12,847\n# analysts across 89 organizations\").\n$anchorText = \"Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n\n$anchorParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $anchorParagraph = $p\n        break\n    }\n}\n\nif ($null -eq $anchorParagraph) {\n    throw \"Could not locate the 'Platform impact' achievement bullet paragraph.\"\n}\n\n# The four bullets to add, in document order.\n$bullet = [char]0x2022\n$pm = [char]0x00B1\n$newBullets = @(\n    \"$bullet Real-time collaboration at national scale\",\n    \"$bullet Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ${pm}4.2% to ${pm}2.1%\",\n    \"$bullet Increased voter turnout prediction accuracy from 71% to 87%\",\n    \"$bullet Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\"\n)\n\n# Insert each bullet as its own new paragraph, immediately after the anchor,\n# chaining so they appear in the same order as in the diff.\n$anchor = $anchorParagraph\nforeach ($text in $newBullets) {\n    $anchor.Range.InsertParagraphAfter()\n    $newPara = $anchor.Next()\n    $newPara.Range.Text = $text\n    $anchor = $newPara\n}\n"}
